# Rename the spectral-value column headers (row 3) from the old short/lowercase
# labels to the new display labels, across all three repeated header blocks
# (I:M, N:R, S:W).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I3").Value = "PGA"
$ws.Range("J3").Value = "PGV"
$ws.Range("K3").Value = "SA(0.3)"
$ws.Range("L3").Value = "SA(1.0)"
$ws.Range("M3").Value = "SA(3.0)"

$ws.Range("N3").Value = "PGA"
$ws.Range("O3").Value = "PGV"
$ws.Range("P3").Value = "SA(0.3)"
$ws.Range("Q3").Value = "SA(1.0)"
$ws.Range("R3").Value = "SA(3.0)"

$ws.Range("S3").Value = "PGA"
$ws.Range("T3").Value = "PGV"
$ws.Range("U3").Value = "SA(0.3)"
$ws.Range("V3").Value = "SA(1.0)"
$ws.Range("W3").Value = "SA(3.0)"

# Reflect the active-cell selection change recorded in the saved view state.
$ws.Range("I4").Select()
